$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: update selection (whole-sheet select, no explicit active cell) ----
$ws1.Range("A1:XFD1048576").Select() | Out-Null

# ---- Sheet2: validCredentials ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "validCredentials"

# Reuse the existing header / data-row styles from Sheet1 (avoids creating new style entries)
$ws1.Range("A1:B1").Copy() | Out-Null
$ws2.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("A2:B2").Copy() | Out-Null
$ws2.Range("A2:B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "abschallengers"
$ws2.Range("B2").Value = "numpyninja17"

$ws2.Columns.Item(1).ColumnWidth = 18.9
$ws2.Columns.Item(2).ColumnWidth = 20.5

$ws2.Range("A3:B10").EntireRow.Select() | Out-Null

# ---- Sheet3: invalidCredential ----
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "invalidCredential"

$ws1.Range("A1:B1").Copy() | Out-Null
$ws3.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("A2:B2").Copy() | Out-Null
$ws3.Range("A2:B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws3.Range("A1").Value = "username"
$ws3.Range("B1").Value = "password"
$ws3.Range("A2").Value = "abschallengers123"
$ws3.Range("B2").Value = "numpyninja17123"

$ws3.Columns.Item(1).ColumnWidth = 18.9
$ws3.Columns.Item(2).ColumnWidth = 20.5

$ws3.Range("A3:B6").EntireRow.Select() | Out-Null
$ws3.Activate() | Out-Null
